# DatasetFitsWorkbook.xlsx update
# - P8/Q8 text tweak ("...on-off..." -> "...on-off-twice...")
# - J10 gets a new note ("use plots from 9d fits")
# - The orange "ready to be fit" status cells in row 9-11 (H,I,J,K,S,T) are
#   cleared back to the green "not yet" fill (G is also extended to match)
# - Row 12 (HighKick) picks up the same green-fill status columns as the
#   other datasets, left blank (not ready yet)
# - Selection / frozen-pane scroll position + a couple of column widths
#   were nudged by the author while reviewing the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5296274   # RGB(146,208,80) / 92D050 - matches existing style "1"
$orange = 49407    # RGB(255,192,0) / FFC000 - matches existing style "2"

# --- text edits -----------------------------------------------------------
$ws.Range("P8").Value = "maybe just do on-off-twice with this"
$ws.Range("Q8").Value = "maybe just do on-off-twice with this"

$ws.Range("J10").Value = "use plots from 9d fits"

# --- row 9 (60h): clear "ready to be fit" status, repaint green -----------
foreach ($col in @("G9","H9","I9","J9","K9","S9","T9")) {
    $ws.Range($col).Value = ""
    $ws.Range($col).Interior.Color = $green
}

# --- row 10 (9d): same, but J10 keeps its new text -------------------------
foreach ($col in @("G10","H10","I10","K10","S10","T10")) {
    $ws.Range($col).Value = ""
    $ws.Range($col).Interior.Color = $green
}
$ws.Range("J10").Interior.Color = $green

# --- row 11 (Endgame): same treatment --------------------------------------
foreach ($col in @("G11","H11","I11","J11","K11","S11","T11")) {
    $ws.Range($col).Value = ""
    $ws.Range($col).Interior.Color = $green
}

# --- row 12 (HighKick): add the orange "not ready" block -------------------
foreach ($col in @("B12","C12","H12","I12","J12","K12","S12","T12","U12","V12")) {
    $ws.Range($col).Value = ""
    $ws.Range($col).Interior.Color = $orange
}

# --- a couple of column widths were widened ---------------------------------
$ws.Columns.Item(16).ColumnWidth = 36
$ws.Columns.Item(17).ColumnWidth = 33.666666666666664

# --- view state: unfreeze scroll back to column B, reselect C31 ------------
[void]$ws.Range("C31").Select()
